$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.4474447272637008
$ws.Cells.Item(2, 3).Value = 0.2054119868599287
$ws.Cells.Item(2, 4).Value = 0.05159773798830969
$ws.Cells.Item(2, 6).Value = 3.681058352749574
$ws.Cells.Item(2, 7).Value = 0.002583181707623449
$ws.Cells.Item(2, 9).Value = 2.0245336421671
$ws.Cells.Item(2, 10).Value = 0.2839615490135117
$ws.Cells.Item(2, 11).Value = 0.6516171022354342
$ws.Cells.Item(2, 13).Value = 0.3363450314127903

$ws.Cells.Item(3, 2).Value = 0.4276025166107331
$ws.Cells.Item(3, 3).Value = 0.1979779558425321
$ws.Cells.Item(3, 4).Value = 0.05207506493093206
$ws.Cells.Item(3, 6).Value = 3.630101729478824
$ws.Cells.Item(3, 7).Value = 0.002587446606767347
$ws.Cells.Item(3, 9).Value = 1.996558129882715
$ws.Cells.Item(3, 10).Value = 0.2810608648078343
$ws.Cells.Item(3, 11).Value = 0.6245140108728151
$ws.Cells.Item(3, 13).Value = 0.3287900991370947

$ws.Cells.Item(4, 2).Value = 0.4157915159369452
$ws.Cells.Item(4, 3).Value = 0.1935475632507604
$ws.Cells.Item(4, 4).Value = 0.05241725743582748
$ws.Cells.Item(4, 6).Value = 3.599933670608294
$ws.Cells.Item(4, 7).Value = 0.002590203233920896
$ws.Cells.Item(4, 9).Value = 1.97986750696586
$ws.Cells.Item(4, 10).Value = 0.2794003651160324
$ws.Cells.Item(4, 11).Value = 0.6083784337099871
$ws.Cells.Item(4, 13).Value = 0.3243874588267914

$ws.Cells.Item(5, 2).Value = 0.4110720037350575
$ws.Cells.Item(5, 3).Value = 0.1917757750057234
$ws.Cells.Item(5, 4).Value = 0.05256908644040692
$ws.Cells.Item(5, 6).Value = 3.587920868856898
$ws.Cells.Item(5, 7).Value = 0.002591361386652643
$ws.Cells.Item(5, 9).Value = 1.973187861491454
$ws.Cells.Item(5, 10).Value = 0.2787539838907946
$ws.Cells.Item(5, 11).Value = 0.6019301188730424
$ws.Cells.Item(5, 13).Value = 0.3226527109224406

$ws.Cells.Item(6, 2).Value = 0.4102939829156753
$ws.Cells.Item(6, 3).Value = 0.1914835994882367
$ws.Cells.Item(6, 4).Value = 0.05259504640304513
$ws.Cells.Item(6, 6).Value = 3.585943102698536
$ws.Cells.Item(6, 7).Value = 0.002591555802550993
$ws.Cells.Item(6, 9).Value = 1.972086063287705
$ws.Cells.Item(6, 10).Value = 0.2786484813275933
$ws.Cells.Item(6, 11).Value = 0.6008670518028509
$ws.Cells.Item(6, 13).Value = 0.3223682425516863

$ws.Cells.Item(7, 2).Value = 0.4157274881293631
$ws.Cells.Item(7, 3).Value = 0.1935235322436881
$ws.Cells.Item(7, 4).Value = 0.052419254872202
$ws.Cells.Item(7, 6).Value = 3.599770525089127
$ws.Cells.Item(7, 7).Value = 0.002590218711958276
$ws.Cells.Item(7, 9).Value = 1.979776929832866
$ws.Cells.Item(7, 10).Value = 0.2793915251947965
$ws.Cells.Item(7, 11).Value = 0.6082909551734019
$ws.Cells.Item(7, 13).Value = 0.3243638230598123

$ws.Cells.Item(8, 2).Value = 0.4405258586890568
$ws.Cells.Item(8, 3).Value = 0.2028208326328809
$ws.Cells.Item(8, 4).Value = 0.05175214683380602
$ws.Cells.Item(8, 6).Value = 3.663255780807063
$ws.Cells.Item(8, 7).Value = 0.00258462367925036
$ws.Cells.Item(8, 9).Value = 2.014786467963532
$ws.Cells.Item(8, 10).Value = 0.2829363611599476
$ws.Cells.Item(8, 11).Value = 0.6421668793763047
$ws.Cells.Item(8, 13).Value = 0.3336910739229353

$ws.Cells.Item(9, 2).Value = 0.4921134292041813
$ws.Cells.Item(9, 3).Value = 0.2221227114975761
$ws.Cells.Item(9, 4).Value = 0.0508322447192171
$ws.Cells.Item(9, 6).Value = 3.796668332799101
$ws.Cells.Item(9, 7).Value = 0.002574741317243445
$ws.Cells.Item(9, 9).Value = 2.087322557746859
$ws.Cells.Item(9, 10).Value = 0.2908458624340682
$ws.Cells.Item(9, 11).Value = 0.7126218557368702
$ws.Cells.Item(9, 13).Value = 0.3538574008304423

$ws.Cells.Item(10, 2).Value = 0.5318301109661832
$ws.Cells.Item(10, 3).Value = 0.2369656842307961
$ws.Cells.Item(10, 4).Value = 0.05039135316724952
$ws.Cells.Item(10, 6).Value = 3.900185158826559
$ws.Cells.Item(10, 7).Value = 0.002568137652582662
$ws.Cells.Item(10, 9).Value = 2.143019548570635
$ws.Cells.Item(10, 10).Value = 0.2972444249421642
$ws.Cells.Item(10, 11).Value = 0.7668606505894218
$ws.Cells.Item(10, 13).Value = 0.3698223649033707

$ws.Cells.Item(11, 2).Value = 0.5502954322676032
$ws.Cells.Item(11, 3).Value = 0.2438640794210016
$ws.Cells.Item(11, 4).Value = 0.0502414423970734
$ws.Cells.Item(11, 6).Value = 3.948485821405399
$ws.Cells.Item(11, 7).Value = 0.002565274566102589
$ws.Cells.Item(11, 9).Value = 2.16888827265791
$ws.Cells.Item(11, 10).Value = 0.3002836596109972
$ws.Cells.Item(11, 11).Value = 0.7920781459226873
$ws.Cells.Item(11, 13).Value = 0.3773360275259847

$ws.Cells.Item(12, 2).Value = 0.5573451469054476
$ws.Cells.Item(12, 3).Value = 0.2464974983916193
$ws.Cells.Item(12, 4).Value = 0.05019192888180868
$ws.Cells.Item(12, 6).Value = 3.966950957744302
$ws.Cells.Item(12, 7).Value = 0.002564210540529839
$ws.Cells.Item(12, 9).Value = 2.17876112527108
$ws.Cells.Item(12, 10).Value = 0.3014530633421657
$ws.Cells.Item(12, 11).Value = 0.801705881211177
$ws.Cells.Item(12, 13).Value = 0.3802174333785047

$ws.Cells.Item(13, 2).Value = 0.5558243152576097
$ws.Cells.Item(13, 3).Value = 0.245929402141968
$ws.Cells.Item(13, 4).Value = 0.05020227032292013
$ws.Cells.Item(13, 6).Value = 3.962966378483117
$ws.Cells.Item(13, 7).Value = 0.002564438802481154
$ws.Cells.Item(13, 9).Value = 2.17663140131252
$ws.Cells.Item(13, 10).Value = 0.3012003874710842
$ws.Cells.Item(13, 11).Value = 0.7996288854006934
$ws.Cells.Item(13, 13).Value = 0.3795952629459336

$ws.Cells.Item(14, 2).Value = 0.5508742677937732
$ws.Cells.Item(14, 3).Value = 0.2440803079091722
$ws.Cells.Item(14, 4).Value = 0.05023722366693306
$ws.Cells.Item(14, 6).Value = 3.950001453652362
$ws.Cells.Item(14, 7).Value = 0.002565186624490575
$ws.Cells.Item(14, 9).Value = 2.169698973877857
$ws.Cells.Item(14, 10).Value = 0.3003794959637816
$ws.Cells.Item(14, 11).Value = 0.7928686527794753
$ws.Cells.Item(14, 13).Value = 0.3775723579166836

$ws.Cells.Item(15, 2).Value = 0.54784968476676
$ws.Cells.Item(15, 3).Value = 0.2429504416552675
$ws.Cells.Item(15, 4).Value = 0.05025957747771059
$ws.Cells.Item(15, 6).Value = 3.942082837272068
$ws.Cells.Item(15, 7).Value = 0.002565647309736527
$ws.Cells.Item(15, 9).Value = 2.165462693442748
$ws.Cells.Item(15, 10).Value = 0.2998790879814095
$ws.Cells.Item(15, 11).Value = 0.788738033136525
$ws.Cells.Item(15, 13).Value = 0.3763379783684968

$ws.Cells.Item(16, 2).Value = 0.5306313731281875
$ws.Cells.Item(16, 3).Value = 0.2365178106769008
$ws.Cells.Item(16, 4).Value = 0.05040216693272725
$ws.Cells.Item(16, 6).Value = 3.89705300988112
$ws.Cells.Item(16, 7).Value = 0.00256832758944725
$ws.Cells.Item(16, 9).Value = 2.141339707378506
$ws.Cells.Item(16, 10).Value = 0.2970483923615888
$ws.Cells.Item(16, 11).Value = 0.7652235907019644
$ws.Cells.Item(16, 13).Value = 0.3693363854114722

$ws.Cells.Item(17, 2).Value = 0.5201704798021183
$ws.Cells.Item(17, 3).Value = 0.2326091401508563
$ws.Cells.Item(17, 4).Value = 0.05050259332389828
$ws.Cells.Item(17, 6).Value = 3.869739146636533
$ws.Cells.Item(17, 7).Value = 0.002570007880349931
$ws.Cells.Item(17, 9).Value = 2.126677580176036
$ws.Cells.Item(17, 10).Value = 0.2953447871216213
$ws.Cells.Item(17, 11).Value = 0.7509377056605047
$ws.Cells.Item(17, 13).Value = 0.3651054761952821

$ws.Cells.Item(18, 2).Value = 0.5141911083198067
$ws.Cells.Item(18, 3).Value = 0.2303747351743084
$ws.Cells.Item(18, 4).Value = 0.05056512535892921
$ws.Cells.Item(18, 6).Value = 3.854142816823554
$ws.Cells.Item(18, 7).Value = 0.00257098761285443
$ws.Cells.Item(18, 9).Value = 2.118294349155391
$ws.Cells.Item(18, 10).Value = 0.2943770110614139
$ws.Cells.Item(18, 11).Value = 0.7427720339358359
$ws.Cells.Item(18, 13).Value = 0.3626956035438766

$ws.Cells.Item(19, 2).Value = 0.5121730273597223
$ws.Cells.Item(19, 3).Value = 0.2296205636693003
$ws.Cells.Item(19, 4).Value = 0.05058711768845114
$ws.Cells.Item(19, 6).Value = 3.848881706287358
$ws.Cells.Item(19, 7).Value = 0.002571321616441372
$ws.Cells.Item(19, 9).Value = 2.115464511532906
$ws.Cells.Item(19, 10).Value = 0.2940514146169022
$ws.Cells.Item(19, 11).Value = 0.7400160642895344
$ws.Cells.Item(19, 13).Value = 0.3618837203145659

$ws.Cells.Item(20, 2).Value = 0.5212801825744577
$ws.Cells.Item(20, 3).Value = 0.2330237999761664
$ws.Cells.Item(20, 4).Value = 0.05049140936084484
$ws.Cells.Item(20, 6).Value = 3.872634962811361
$ws.Cells.Item(20, 7).Value = 0.002569827637820377
$ws.Cells.Item(20, 9).Value = 2.12823320684609
$ws.Cells.Item(20, 10).Value = 0.2955248870251523
$ws.Cells.Item(20, 11).Value = 0.7524531630306797
$ws.Cells.Item(20, 13).Value = 0.3655534173951125

$ws.Cells.Item(21, 2).Value = 0.5523266611137103
$ws.Cells.Item(21, 3).Value = 0.2446228568682614
$ws.Cells.Item(21, 4).Value = 0.05022676037148699
$ws.Cells.Item(21, 6).Value = 3.95380481895603
$ws.Cells.Item(21, 7).Value = 0.002564966424735031
$ws.Cells.Item(21, 9).Value = 2.171733104500078
$ws.Cells.Item(21, 10).Value = 0.3006201090257292
$ws.Cells.Item(21, 11).Value = 0.7948521671884805
$ws.Cells.Item(21, 13).Value = 0.3781655529118524

$ws.Cells.Item(22, 2).Value = 0.5729513422332673
$ws.Cells.Item(22, 3).Value = 0.2523268238566061
$ws.Cells.Item(22, 4).Value = 0.05009606807091416
$ws.Cells.Item(22, 6).Value = 4.007872755525653
$ws.Cells.Item(22, 7).Value = 0.002561906820565693
$ws.Cells.Item(22, 9).Value = 2.200611406460723
$ws.Cells.Item(22, 10).Value = 0.3040580468916403
$ws.Cells.Item(22, 11).Value = 0.8230195368580553
$ws.Cells.Item(22, 13).Value = 0.3866190105111116

$ws.Cells.Item(23, 2).Value = 0.5619129762770285
$ws.Cells.Item(23, 3).Value = 0.2482037497987051
$ws.Cells.Item(23, 4).Value = 0.05016196279397889
$ws.Cells.Item(23, 6).Value = 3.978922239153604
$ws.Cells.Item(23, 7).Value = 0.002563529073973532
$ws.Cells.Item(23, 9).Value = 2.185157318645238
$ws.Cells.Item(23, 10).Value = 0.302213269000589
$ws.Cells.Item(23, 11).Value = 0.8079441852052582
$ws.Cells.Item(23, 13).Value = 0.3820879501964569

$ws.Cells.Item(24, 2).Value = 0.5207783777603083
$ws.Cells.Item(24, 3).Value = 0.2328362925892975
$ws.Cells.Item(24, 4).Value = 0.05049645069085784
$ws.Cells.Item(24, 6).Value = 3.871325431772249
$ws.Cells.Item(24, 7).Value = 0.002569909082649377
$ws.Cells.Item(24, 9).Value = 2.127529764081871
$ws.Cells.Item(24, 10).Value = 0.2954434275763589
$ws.Cells.Item(24, 11).Value = 0.7517678768952294
$ws.Cells.Item(24, 13).Value = 0.3653508329937551

$ws.Cells.Item(25, 2).Value = 0.4778397281536968
$ws.Cells.Item(25, 3).Value = 0.2167855129421525
$ws.Cells.Item(25, 4).Value = 0.05103971610161295
$ws.Cells.Item(25, 6).Value = 3.759615580763949
$ws.Cells.Item(25, 7).Value = 0.002577298875631157
$ws.Cells.Item(25, 9).Value = 2.067280069950442
$ws.Cells.Item(25, 10).Value = 0.2886032172902446
$ws.Cells.Item(25, 11).Value = 0.6931288193439968
$ws.Cells.Item(25, 13).Value = 0.3482005362370515
